{"js": "const pairs = [\n  [\"81\u00f79=\", \"48\u00f74=\"],\n  [\"12\u00f76=\", \"56\u00f73=\"],\n  [\"36\u00f79=\", \"88\u00f74=\"],\n  [\"68\u00f73=\", \"65\u00f79=\"],\n  [\"56\u00f79=\", \"27\u00f76=\"],\n  [\"99\u00f76=\", \"86\u00f73=\"],\n  [\"86\u00f76=\", \"49\u00f79=\"],\n  [\"67\u00f73=\", \"29\u00f76=\"],\n  [\"42\u00f74=\", \"94\u00f75=\"],\n  [\"17\u00f78=\", \"30\u00f77=\"],\n  [\"33\u00f72=\", \"93\u00f74=\"],\n  [\"24\u00f76=\", \"92\u00f74=\"],\n  [\"60\u00f76=\", \"72\u00f73=\"],\n  [\"11\u00f74=\", \"97\u00f79=\"],\n  [\"33\u00f74=\", \"37\u00f77=\"],\n  [\"95\u00f75=\", \"12\u00f74=\"],\n  [\"56\u00f76=\", \"69\u00f79=\"],\n  [\"18\u00f77=\", \"84\u00f77=\"],\n  [\"42\u00f76=\", \"91\u00f79=\"],\n  [\"43\u00f77=\", \"16\u00f73=\"],\n  [\"39\u00f78=\", \"19\u00f74=\"],\n  [\"45\u00f76=\", \"36\u00f76=\"],\n  [\"44\u00f78=\", \"44\u00f75=\"],\n  [\"17\u00f75=\", \"39\u00f73=\"],\n  [\"35\u00f74=\", \"58\u00f75=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-expression text in each table cell of the\n# \"two-digit number divided by one-digit number\" worksheet.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"81\u00f79=\", \"48\u00f74=\"),\n    @(\"12\u00f76=\", \"56\u00f73=\"),\n    @(\"36\u00f79=\", \"88\u00f74=\"),\n    @(\"68\u00f73=\", \"65\u00f79=\"),\n    @(\"56\u00f79=\", \"27\u00f76=\"),\n    @(\"99\u00f76=\", \"86\u00f73=\"),\n    @(\"86\u00f76=\", \"49\u00f79=\"),\n    @(\"67\u00f73=\", \"29\u00f76=\"),\n    @(\"42\u00f74=\", \"94\u00f75=\"),\n    @(\"17\u00f78=\", \"30\u00f77=\"),\n    @(\"33\u00f72=\", \"93\u00f74=\"),\n    @(\"24\u00f76=\", \"92\u00f74=\"),\n    @(\"60\u00f76=\", \"72\u00f73=\"),\n    @(\"11\u00f74=\", \"97\u00f79=\"),\n    @(\"33\u00f74=\", \"37\u00f77=\"),\n    @(\"95\u00f75=\", \"12\u00f74=\"),\n    @(\"56\u00f76=\", \"69\u00f79=\"),\n    @(\"18\u00f77=\", \"84\u00f77=\"),\n    @(\"42\u00f76=\", \"91\u00f79=\"),\n    @(\"43\u00f77=\", \"16\u00f73=\"),\n    @(\"39\u00f78=\", \"19\u00f74=\"),\n    @(\"45\u00f76=\", \"36\u00f76=\"),\n    @(\"44\u00f78=\", \"44\u00f75=\"),\n    @(\"17\u00f75=\", \"39\u00f73=\"),\n    @(\"35\u00f74=\", \"58\u00f75=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n"}
